# Generate Report for Handoff
#
# The localization run finished: zh-cn/de-de status flips from
# "In Translation" -> "Ready for handoff", and the associated timestamps
# advance a few seconds (new xliff generation / handoff times).
#
# Sheets:
#   Overview  - E2 (zh-cn status), F2 (de-de status), G2 (latest HO xliff
#               generate date)
#   zh-cn     - C2 (status), H2 (latest handoff datetime)
#   de-de     - C2 (status), H2 (latest handoff datetime)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-09-01 17:08:37"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-09-01 17:08:32"

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-09-01 17:08:37"

# The status column got wider now that "Ready for handoff" (18 chars) is
# longer than "In Translation" (14 chars) -- match the report generator's
# auto-fit by nudging the column width (closest value this host's
# character-width quantization can reproduce).
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
